$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 159.11111
$ws.Range("I38").Value = 116.5
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 349.5
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 22.5
$ws.Range("N38").Value = -2244

# Row 58
$ws.Range("H58").Value = 2020.5714
$ws.Range("I58").Value = 440.66666
$ws.Range("J58").Value = 11500
$ws.Range("K58").Value = 1321.99998
$ws.Range("L58").Value = 34500
$ws.Range("M58").Value = -1171.99998
$ws.Range("N58").Value = -34800

# Row 64
$ws.Range("H64").Value = 3320.4
$ws.Range("I64").Value = 2993.7144
$ws.Range("J64").Value = 3606.25
$ws.Range("K64").Value = 2993.7144
$ws.Range("L64").Value = 3606.25
$ws.Range("M64").Value = -2745.7144
$ws.Range("N64").Value = -4102.25

# Row 67
$ws.Range("H67").Value = 3320.4
$ws.Range("I67").Value = 2993.7144
$ws.Range("J67").Value = 3606.25
$ws.Range("K67").Value = 2993.7144
$ws.Range("L67").Value = 3606.25
$ws.Range("M67").Value = -2135.7144
$ws.Range("N67").Value = -5322.25

# Row 94
$ws.Range("H94").Value = 1122.25
$ws.Range("I94").Value = 1122.25
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1122.25
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = -671.25

# Row 106
$ws.Range("H106").Value = 4570.4736
$ws.Range("I106").Value = 1284.6471
$ws.Range("K106").Value = 1284.6471
$ws.Range("M106").Value = -653.6470999999999

# Row 116
$ws.Range("H116").Value = 20164638
$ws.Range("I116").Value = 17939858
$ws.Range("J116").Value = 22241098
$ws.Range("K116").Value = 17939858
$ws.Range("L116").Value = 22241098
$ws.Range("M116").Value = -17936416
$ws.Range("N116").Value = -22247982

# Row 132
$ws.Range("H132").Value = 5190.0605
$ws.Range("I132").Value = 4661.483
$ws.Range("J132").Value = 9022.25
$ws.Range("K132").Value = 13984.449
$ws.Range("L132").Value = 27066.75
$ws.Range("M132").Value = -11454.449
$ws.Range("N132").Value = -32126.75

# Row 137
$ws.Range("H137").Value = 2323.6843
$ws.Range("I137").Value = 2126.5293
$ws.Range("K137").Value = 6379.5879
$ws.Range("M137").Value = -3829.5879

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1374.3334
$ws.Range("I2").Value = 1357.4166
$ws.Range("K2").Value = 1357.4166
$ws.Range("M2").Value = -1244.4166

# Row 32
$ws.Range("H32").Value = 2635.9434
$ws.Range("I32").Value = 2584.2
$ws.Range("K32").Value = 2584.2
$ws.Range("M32").Value = -2297.2

# Row 61
$ws.Range("H61").Value = 66669652
$ws.Range("I61").Value = 111113944
$ws.Range("K61").Value = 111113944
$ws.Range("M61").Value = -111113732

# Row 81
$ws.Range("H81").Value = 90000
$ws.Range("J81").Value = 90000
$ws.Range("L81").Value = 90000
$ws.Range("N81").Value = -91996

# Row 84
$ws.Range("H84").Value = 90000
$ws.Range("J84").Value = 90000
$ws.Range("L84").Value = 270000
$ws.Range("N84").Value = -279984

# Row 86
$ws.Range("H86").Value = 999990
$ws.Range("J86").Value = 999990
$ws.Range("L86").Value = 999990
$ws.Range("N86").Value = -1002362

# Row 88
$ws.Range("H88").Value = 5953668.5
$ws.Range("I88").Value = 12821030
$ws.Range("K88").Value = 12821030
$ws.Range("M88").Value = -12820624

# Row 89
$ws.Range("H89").Value = 999990
$ws.Range("J89").Value = 999990
$ws.Range("L89").Value = 2999970
$ws.Range("N89").Value = -3011826

# Row 91
$ws.Range("H91").Value = 5953668.5
$ws.Range("I91").Value = 12821030
$ws.Range("K91").Value = 12821030
$ws.Range("M91").Value = -12819626

# Row 116
$ws.Range("H116").Value = 1374.3334
$ws.Range("I116").Value = 1357.4166
$ws.Range("K116").Value = 1357.4166
$ws.Range("M116").Value = 936.5834

# Row 122
$ws.Range("H122").Value = 14495725
$ws.Range("I122").Value = 17546404
$ws.Range("K122").Value = 52639212
$ws.Range("M122").Value = -52636762

# Row 136
$ws.Range("H136").Value = 66669652
$ws.Range("I136").Value = 111113944
$ws.Range("K136").Value = 333341832
$ws.Range("M136").Value = -333339282

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1374.3334
$ws.Range("I3").Value = 1357.4166
$ws.Range("K3").Value = 1357.4166
$ws.Range("M3").Value = -1243.4166

# Row 86
$ws.Range("H86").Value = 2353.75
$ws.Range("I86").Value = 2163.8096
$ws.Range("J86").Value = 2716.3635
$ws.Range("K86").Value = 2163.8096
$ws.Range("L86").Value = 2716.3635
$ws.Range("M86").Value = -1040.8096
$ws.Range("N86").Value = -4962.363499999999

# Row 89
$ws.Range("H89").Value = 2353.75
$ws.Range("I89").Value = 2163.8096
$ws.Range("J89").Value = 2716.3635
$ws.Range("K89").Value = 10819.048
$ws.Range("L89").Value = 13581.8175
$ws.Range("M89").Value = -5203.048000000001
$ws.Range("N89").Value = -24813.8175

# Row 107
$ws.Range("H107").Value = 15152468
$ws.Range("J107").Value = 125000510
$ws.Range("L107").Value = 125000510
$ws.Range("N107").Value = -125004350

# Row 134
$ws.Range("H134").Value = 1409.4762
$ws.Range("I134").Value = 1329.95
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3989.85
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1454.85
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1343.2222
$ws.Range("I16").Value = 1082.5883
$ws.Range("K16").Value = 1082.5883
$ws.Range("M16").Value = -795.5882999999999

# Row 31
$ws.Range("H31").Value = 3102.2856
$ws.Range("I31").Value = 2543.5
$ws.Range("K31").Value = 2543.5
$ws.Range("M31").Value = -2248.5

# Row 34
$ws.Range("H34").Value = 3102.2856
$ws.Range("I34").Value = 2543.5
$ws.Range("K34").Value = 2543.5
$ws.Range("M34").Value = -2341.5

# Row 68
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498

# Row 71
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488

# Row 86
$ws.Range("H86").Value = 74349.586
$ws.Range("I86").Value = 108562
$ws.Range("J86").Value = 5924.75
$ws.Range("K86").Value = 108562
$ws.Range("L86").Value = 5924.75
$ws.Range("M86").Value = -107439
$ws.Range("N86").Value = -8170.75

# Row 89
$ws.Range("H89").Value = 74349.586
$ws.Range("I89").Value = 108562
$ws.Range("J89").Value = 5924.75
$ws.Range("K89").Value = 542810
$ws.Range("L89").Value = 29623.75
$ws.Range("M89").Value = -537194
$ws.Range("N89").Value = -40855.75

# Row 113
$ws.Range("H113").Value = 1343.2222
$ws.Range("I113").Value = 1082.5883
$ws.Range("K113").Value = 1082.5883
$ws.Range("M113").Value = 1087.4117

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = $null
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = 0

# Row 134
$ws.Range("H134").Value = 2830.5
$ws.Range("I134").Value = 2542.4546
$ws.Range("K134").Value = 7627.3638
$ws.Range("M134").Value = -5092.3638

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 1456.2
$ws.Range("J97").Value = 1387.1666
$ws.Range("L97").Value = 4161.4998
$ws.Range("N97").Value = -5153.4998

$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 7521.143
$ws.Range("I41").Value = 5474.75
$ws.Range("J41").Value = 10249.667
$ws.Range("K41").Value = 5474.75
$ws.Range("L41").Value = 10249.667
$ws.Range("M41").Value = -5119.75
$ws.Range("N41").Value = -10959.667

# Row 113
$ws.Range("H113").Value = 3043.3076
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 4351.625
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 4351.625
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -8691.625

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 25003582
$ws.Range("I7").Value = 41669532
$ws.Range("K7").Value = 41669532
$ws.Range("M7").Value = -41669420

# Row 22
$ws.Range("H22").Value = 2688575.2
$ws.Range("I22").Value = 387
$ws.Range("J22").Value = 3584638
$ws.Range("K22").Value = 387
$ws.Range("L22").Value = 3584638
$ws.Range("M22").Value = -92
$ws.Range("N22").Value = -3585228

# Row 27
$ws.Range("H27").Value = 2688575.2
$ws.Range("I27").Value = 387
$ws.Range("J27").Value = 3584638
$ws.Range("K27").Value = 387
$ws.Range("L27").Value = 3584638
$ws.Range("M27").Value = -280
$ws.Range("N27").Value = -3584852

# Row 40
$ws.Range("H40").Value = 2911.9
$ws.Range("I40").Value = 2481.2856
$ws.Range("J40").Value = 3916.6667
$ws.Range("K40").Value = 2481.2856
$ws.Range("L40").Value = 3916.6667
$ws.Range("M40").Value = -2345.2856
$ws.Range("N40").Value = -4188.6667

# Row 122
$ws.Range("H122").Value = 3514.0312
$ws.Range("I122").Value = 2834.5
$ws.Range("K122").Value = 8503.5
$ws.Range("M122").Value = -6053.5

# Row 126
$ws.Range("H126").Value = 25003582
$ws.Range("I126").Value = 41669532
$ws.Range("K126").Value = 125008596
$ws.Range("M126").Value = -125006126

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1664.1
$ws.Range("I126").Value = 1493.3125
$ws.Range("J126").Value = 2347.25
$ws.Range("K126").Value = 4479.9375
$ws.Range("L126").Value = 7041.75
$ws.Range("M126").Value = -2009.9375
$ws.Range("N126").Value = -11981.75

# Row 132
$ws.Range("H132").Value = 6999.6924
$ws.Range("I132").Value = 7624.625
$ws.Range("K132").Value = 22873.875
$ws.Range("M132").Value = -20343.875

# Row 136
$ws.Range("H136").Value = 4885.4614
$ws.Range("I136").Value = 2073.1428
$ws.Range("K136").Value = 6219.428400000001
$ws.Range("M136").Value = -3669.428400000001
